# The "c_Name" column (column B: id | c_Name | f_Name | values) is a
# constant column (all rows = "base") that is no longer needed; remove it
# and shift f_Name/values left, matching the data reduction described in
# the commit (simplifying the input data fixture used for model-result
# comparison tests).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("B").Delete()

# Update the recorded selection to match the edited file (C15).
$ws.Range("C15").Select()
